# Expense Tracker - "Create Expense APIs": append new Income entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Income")

# Find the first empty row below the existing data (rows 1-3 are used).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$newEntries = @(
    @{ Source = "Salary"; Amount = 4000; Date = 45728.229537037034 },
    @{ Source = "Salary"; Amount = 5000; Date = 45728.229537037034 }
)

foreach ($entry in $newEntries) {
    $row = $lastRow + 1

    $ws.Cells.Item($row, 1).Value2 = $entry.Source
    $ws.Cells.Item($row, 2).Value2 = $entry.Amount

    # Write the date, then copy the existing date column's formatting so the
    # new cell keeps the same date number format as the rows above it.
    $ws.Cells.Item($row, 3).Value2 = $entry.Date
    $ws.Range("C3").Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
    $ws.Cells.Item($row, 3).Value2 = $entry.Date

    $lastRow = $row
}

$excel.CutCopyMode = $false
